# [imp] update style of btn GPS-23
#
# This reproduces the target edit to Sheet1:
#  - Removes the two extra "button" picture copies (Picture 4 / Picture 5) that were
#    anchored on rows 2 and 3, keeping only the original picture (Picture 3) on row 1.
#  - Updates row 1: D1 changes from 37 to 168 (A1 keeps its "Apple" text).
#  - Clears out all of the leftover row 2 / row 3 data (they become essentially blank
#    rows, with only the still-styled-but-empty date cell left behind in column G).
#
# Shared-string table cleanup (dropping the now-unused "Pizza"/"orange" entries) and
# shared-string index renumbering happen automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drawings: remove the two duplicated button pictures, keep the first one ---
$ws.Shapes.Item("Picture 5").Delete()
$ws.Shapes.Item("Picture 4").Delete()

# --- Row 1: bump the D1 number from 37 to 168 ---
$ws.Range("D1").Value = 168

# --- Row 2: wipe out all the old row data, leaving only the styled-but-empty G2 cell ---
$ws.Range("A2:E2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()

# --- Row 3: same cleanup as row 2 ---
$ws.Range("A3:E3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
